# Apply the "RIP funcionando entre Noroeste & Centro" edit:
# On the VLAN sheet, duplicate the Sub-Interfaz/IP Network/Mascara block
# (columns D:F) from the first mini-table (rows 2-6) into the second
# mini-table (rows 10-14), with the new g0/0.30-.95 sub-interfaces and
# 130.45.128.0/26-based addressing for the Noroeste<->Centro link.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VLAN")

# --- Copy the cell formatting from the first table (rows 2-6) onto the
#     new D:F cells of the second table (rows 10-14), row by row. ---
$ws.Range("D2:F2").Copy() | Out-Null
$ws.Range("D10:F10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("D3:F3").Copy() | Out-Null
$ws.Range("D11:F11").PasteSpecial(-4122) | Out-Null

$ws.Range("D4:F4").Copy() | Out-Null
$ws.Range("D12:F12").PasteSpecial(-4122) | Out-Null

$ws.Range("D5:F5").Copy() | Out-Null
$ws.Range("D13:F13").PasteSpecial(-4122) | Out-Null

$ws.Range("D6:F6").Copy() | Out-Null
$ws.Range("D14:F14").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Header row for the second table (row 10) ---
$ws.Range("D10").Value = "Sub-Interfaz"
$ws.Range("E10").Value = "IP Network"
$ws.Range("F10").Value = "Mascara"

# --- Data entered column by column (Sub-Interfaz, then IP Network,
#     then Mascara) to match the typed-in-order of the new rows. ---
$ws.Range("D11").Value = "g0/0.30"
$ws.Range("D12").Value = "g0/0.35"
$ws.Range("D13").Value = "g0/0.40"
$ws.Range("D14").Value = "g0/0.95"

$ws.Range("E11").Value = "130.45.128.1/26"
$ws.Range("E12").Value = "130.45.128.65/26"
$ws.Range("E13").Value = "130.45.128.129/26"
$ws.Range("E14").Value = "130.45.128.193/26"

$ws.Range("F11").Value = "255.255.255.192"
$ws.Range("F12").Value = "255.255.255.192"
$ws.Range("F13").Value = "255.255.255.192"
$ws.Range("F14").Value = "255.255.255.192"

# --- Column E widened slightly now the IP strings are a bit longer ---
$ws.Range("E1").EntireColumn.ColumnWidth = 31.140625

# --- Cursor/selection bookkeeping that Excel records on save ---
$ws1 = $wb.Worksheets.Item("Ruteo")
$ws1.Range("C9").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("E12").Select() | Out-Null
